# Update countries & provincias Spain
#
# Refreshes the COVID-19 "paises" data dump: bumps the "last updated"
# timestamp, swaps Burkina Faso/Guinea and Mozambique/Siria into their
# correctly sorted rows (Burkina Faso gets fresh figures, while Guinea,
# Mozambique and Siria keep their previous figures - they simply moved one
# row), and refreshes a handful of other countries' statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 21:22"

# Estados Unidos (row 4): Casos totales, Nuevos casos, Recuperados, Casos criticos, Muertes
$ws.Cells.Item(4, 2).Value = 782987
$ws.Cells.Item(4, 3).Value = 18351
$ws.Cells.Item(4, 5).Value = 669440
$ws.Cells.Item(4, 7).Value = 1202
$ws.Cells.Item(4, 8).Value = 41777

# Brasil (row 15)
$ws.Cells.Item(15, 2).Value = 39681
$ws.Cells.Item(15, 3).Value = 1027
$ws.Cells.Item(15, 5).Value = 15039
$ws.Cells.Item(15, 7).Value = 50
$ws.Cells.Item(15, 8).Value = 2512

# Sudafrica (row 54)
$ws.Cells.Item(54, 2).Value = 3300
$ws.Cells.Item(54, 3).Value = 142
$ws.Cells.Item(54, 4).Value = 1055
$ws.Cells.Item(54, 5).Value = 2187
$ws.Cells.Item(54, 7).Value = 4
$ws.Cells.Item(54, 8).Value = 58

# Nueva Zelanda (row 71)
$ws.Cells.Item(71, 6).Value = 2

# Principado de Andorra (row 91)
$ws.Cells.Item(91, 2).Value = 717
$ws.Cells.Item(91, 3).Value = 4
$ws.Cells.Item(91, 4).Value = 248
$ws.Cells.Item(91, 5).Value = 433

# Costa Rica (row 94)
$ws.Cells.Item(94, 2).Value = 662
$ws.Cells.Item(94, 3).Value = 2
$ws.Cells.Item(94, 4).Value = 124
$ws.Cells.Item(94, 5).Value = 532
$ws.Cells.Item(94, 6).Value = 8

# Row 98 becomes Burkina Faso (fresh data); row 99 becomes Guinea
# (previous Guinea figures, unchanged) - the two countries swap rows in
# the sorted listing.
$ws.Cells.Item(98, 1).Value = "Burkina Faso"
$ws.Cells.Item(98, 2).Value = 581
$ws.Cells.Item(98, 3).Value = 5
$ws.Cells.Item(98, 4).Value = 357
$ws.Cells.Item(98, 5).Value = 186
$ws.Cells.Item(98, 7).Value = 2
$ws.Cells.Item(98, 8).Value = 38

$ws.Cells.Item(99, 1).Value = "Guinea"
$ws.Cells.Item(99, 2).Value = 579
$ws.Cells.Item(99, 3).Value = 0
$ws.Cells.Item(99, 4).Value = 87
$ws.Cells.Item(99, 5).Value = 487
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 5

# Row 168 becomes Mozambique; row 169 becomes Siria (figures swap, values
# themselves unchanged).
$ws.Cells.Item(168, 1).Value = "Mozambique"
$ws.Cells.Item(168, 4).Value = 8
$ws.Cells.Item(168, 8).Value = 0

$ws.Cells.Item(169, 1).Value = "Siria"
$ws.Cells.Item(169, 4).Value = 5
$ws.Cells.Item(169, 8).Value = 3
